# Trade #72 closed at 2026-02-17 12:54:51 - unknown UNKNOWN +0.000%
#
# 1) Summary sheet: bump Total Trades (B6) and recompute Win Rate % (B9).
# 2) Strategy Status sheet: bump MarketMaking Trades (D4) and Win Rate % (G4).
# 3) All Trades / MarketMaking sheets: append new trade row (row 73).

$wb = $excel.ActiveWorkbook

# --- 1) Summary ---------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 72
$summary.Range("B9").Value = 44.44

# --- 2) Strategy Status --------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 72
$status.Range("G4").Value = 44.44

# --- 3) New trade row (appended identically to both trade-log sheets) ---
$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A73").Value = 72
    # Leading apostrophe forces these to stay literal text instead of
    # being auto-recognized/converted into a date serial number.
    $ws.Range("B73").Value = "'2026-02-17"
    $ws.Range("C73").Value = "12:54:44"
    $ws.Range("D73").Value = "MarketMaking"
    $ws.Range("E73").Value = "UP"
    $ws.Range("F73").Value = 0.03
    $ws.Range("G73").Value = 0.03
    $ws.Range("H73").Value = "CLOSED"
    $ws.Range("I73").Value = 0
    $ws.Range("J73").Value = 0
    $ws.Range("K73").Value = 100.21
    $ws.Range("L73").Value = 0
    $ws.Range("M73").Value = 0
    $ws.Range("N73").Value = 0.6
    $ws.Range("O73").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P73").Value = "early_exit"
    $ws.Range("Q73").Value = 0.11
}
